$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the existing header
# formatting (bold, centered, bordered) used by the other header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I ("I0") and J ("IF")
$i0 = @(1, 7, 8, 1, 1, 1, 1, 1, 1, 1, 1, 5, 1, 1, 1)
$if = @(5, 7, 9, 4, 4, 6, 6, 5, 7, 6, 7, 7, 2, 3, 2)

for ($k = 0; $k -lt 15; $k++) {
    $row = $k + 2
    $ws.Cells.Item($row, 9).Value = $i0[$k]
    $ws.Cells.Item($row, 10).Value = $if[$k]
}
